# Updating Gym workouts for 10/01/2018
# Appends 11 new workout-log rows (342..352 / sheet rows 343..353) for the
# 10 Jan 2018 (DateId 40) session, introducing one new exercise name
# ("Wall Squats") along the way.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: ExerciseId, DateId, ExerciseDate(serial), ExerciseWeek,
#             ExerciseMonth, ExerciseYear, ExerciseDay, ExerciseName,
#             Weight, Sets, Reps, TrainingArea
$newRows = @(
    @(342, 40, 43110, 2, "January", 2018, "Wednesday", "Dumbell Step Up",      32,  4, 12, "Legs"),
    @(343, 40, 43110, 2, "January", 2018, "Wednesday", "Bodyweight Dip",      102,  4,  8, "Arms"),
    @(344, 40, 43110, 2, "January", 2018, "Wednesday", "Bodyweight Pull-up",  102,  4,  5, "Arms"),
    @(345, 40, 43110, 2, "January", 2018, "Wednesday", "Box jumps",           102,  3, 10, "Legs"),
    @(346, 40, 43110, 2, "January", 2018, "Wednesday", "Press ups",             0,  5, 10, "Chest"),
    @(347, 40, 43110, 2, "January", 2018, "Wednesday", "Squat Snatch",          8,  4, 10, "Legs"),
    @(348, 40, 43110, 2, "January", 2018, "Wednesday", "Wall Squats",           0,  4, 30, "Legs"),
    @(349, 40, 43110, 2, "January", 2018, "Wednesday", "Plank",                 0,  4, 30, "Core"),
    @(350, 40, 43110, 2, "January", 2018, "Wednesday", "Right Plank",           0,  4, 30, "Core"),
    @(351, 40, 43110, 2, "January", 2018, "Wednesday", "Left Plank",            0,  4, 30, "Core"),
    @(352, 40, 43110, 2, "January", 2018, "Wednesday", "Sled Pushes",          20,  5,  4, "Core")
)

$startRow = 343
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $data = $newRows[$i]

    $ws.Cells.Item($r, 1).Value = $data[0]    # A ExerciseId
    $ws.Cells.Item($r, 2).Value = $data[1]    # B DateId
    $ws.Cells.Item($r, 3).Value = $data[2]    # C Exercise Date
    $ws.Cells.Item($r, 4).Value = $data[3]    # D Exercise Week
    $ws.Cells.Item($r, 5).Value = $data[4]    # E Exercise Month
    $ws.Cells.Item($r, 6).Value = $data[5]    # F Exercise Year
    $ws.Cells.Item($r, 7).Value = $data[6]    # G Exercise Day
    $ws.Cells.Item($r, 8).Value = $data[7]    # H Exercise Name
    $ws.Cells.Item($r, 9).Value = $data[8]    # I Weight
    $ws.Cells.Item($r, 10).Value = $data[9]   # J Sets
    $ws.Cells.Item($r, 11).Value = $data[10]  # K Reps
    $ws.Cells.Item($r, 12).Value = $data[11]  # L TrainingArea
}

# Re-anchor the view: keep the header row frozen and land the active
# selection on the last new row's TrainingArea cell, matching where the
# edit left off.
$win = $excel.ActiveWindow
$win.FreezePanes = $false
$null = $ws.Range("A2").Select()
$win.FreezePanes = $true
$null = $ws.Range("M353").Select()
